# Build site at 2022-09-26 16:07:08 UTC
#
# Re-shuffles the "8800009" course-syllabus sheet: a row's worth of content
# that used to sit between "Objetivos:" and "Programa resumido:" is gone,
# several labels/paragraphs slide up one row, a couple of cells get replaced
# outright (professor name shows up twice, a long Portuguese program becomes
# "Semestral", etc.), the old B17:C17 paragraph is dropped, B18:C18 gains the
# professor-name text, and the trailing Bibliografia row disappears entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper behaviour notes (established by experimentation against the
# headless engine):
#   * Range.Value auto-coerces digit-only / date-shaped text into a real
#     number or date serial. Two of the edits below ("01/01/2017") must
#     stay text, so we flip NumberFormat to "@" first, assign the value,
#     then PasteSpecial(xlPasteFormats = -4122) a same-styled neighbour
#     cell back on top to restore the sheet's normal (non-"@") style.
#   * A brand-new cell in a row that previously had no B/C entry doesn't
#     know which style (s="2"/s="3") to inherit, so we copy formats from
#     an existing B/C cell in the same column first.
# ---------------------------------------------------------------------

# --- B10/C10: long Portuguese "Objetivos" paragraph -> professor name ---
$ws.Range("B10").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"
$ws.Range("C10").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"

# --- Row 13: gains an A13 label; B13/C13 text changes ---
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: label + short-syllabus EN text ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "✶Vocal classification. ✶Breathing into the singing. ✶Placing the emission in Bocca Chiusa. ✶Choral singing in unison. ✶Choral singing in canon. ✶Choral singing in other polyphonic formations. ✶Coral reading."
$ws.Range("C14").Value = "✶Vocal classification. ✶Breathing into the singing. ✶Placing the emission in Bocca Chiusa. ✶Choral singing in unison. ✶Choral singing in canon. ✶Choral singing in other polyphonic formations. ✶Coral reading."

# --- Row 15: label "Programa:"; B15/C15 become the date "01/01/2017" (must stay text) ---
$ws.Range("A15").Value = "Programa:"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2017"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2017"
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: label "Syllabus:"; long EN syllabus text stays in place ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "✶Vocal classification.✶Breathing into the singing - Exercises for locating low and average breathing. Support and air column. ✶Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. ✶Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. ✶Choral singing in canon. ✶Choral singing in other polyphonic formations. ✶Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. ✶Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission."
$ws.Range("C16").Value = "✶Vocal classification.✶Breathing into the singing - Exercises for locating low and average breathing. Support and air column. ✶Placing the emission in Bocca Chiusa. - relaxation of the mandible, tongue position, the soft palate suspension, local perception exercises for where the voice is being put, support connection and vocal emission, passing the Bocca Chiusa for vowels and other nasal and guttural sounds. ✶Choral singing in unison. - The choral singing in unison, tuning, timbre uniformity, rhythmic precision. ✶Choral singing in canon. ✶Choral singing in other polyphonic formations. ✶Coral reading.- testing of harmonic relaying and listening to 1st view the various voices, memorization, music theory basics. ✶Assembling and improvement of pieces - promoting the application of learnt techniques. Connection between diaphragm and vocal emission."

# --- Row 17: label becomes "Avaliação:"; drop the old B17/C17 paragraph and its 120pt height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).RowHeight = 15

# --- Row 18: label "Método:"; gains B18/C18 (professor name again) ---
$ws.Range("A18").Value = "Método:"

$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"

$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"

$ws.Rows.Item(18).RowHeight = 60

# --- Rows 19-21: labels shift up one slot (B/C content unaffected) ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Old row 22 (Bibliografia: + giant reference list) is gone entirely ---
$ws.Rows.Item(22).Delete()
